$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.867.94'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '1.585.07'
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.81'
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.478'
$ws.Range("E7").Value = '  -3.78%  '
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.07'
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").Value = '1.806.16'
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = '1.584.54'
$ws.Range("E13").Value = '  -2.22%  '
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("E15").Value = '  -2.89%  '
$ws.Range("D16").Value = '25.845.00'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '59.77'
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("E21").Value = '  -1.66%  '
$ws.Range("E22").Value = '  -1.78%  '
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.09'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.10'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("E29").Value = '  -3.06%  '
$ws.Range("E30").Value = '  -5.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0471'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("D36").Value = '1.101.71'
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("E39").Value = '  -2.40%  '
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.822'
$ws.Range("E41").Value = '  +8.35%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.778'
$ws.Range("E42").Value = '  -7.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.18'
$ws.Range("E43").Value = '  +1.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.80'
$ws.Range("E44").Value = '  -3.99%  '
$ws.Range("D45").Value = '1.719.94'
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D46").Value = '0.0₆0111'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.24'
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.407'
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("E51").Value = '  -0.09%  '
